$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing inputs (B2/B3) to new currency-style values ---
$ws.Range("B2").Value = 1.5501
$ws.Range("B3").Value = 1.5505

# --- Exit ratio / price labels + values (D4:E11) written first so the
#     shared-string table picks these up before "Result" / the G-column
#     helper labels. ---
$ws.Range("D4").Value = "Exit 1 Ratio"
$ws.Range("E4").Value = 0.25

$ws.Range("D5").Value = "Exit 1 Price"
$ws.Range("E5").Value = 1.5494000000000001

$ws.Range("D6").Value = "Exit 2 Ratio"
$ws.Range("E6").Value = 0.75

$ws.Range("D7").Value = "Exit 2 Price"
$ws.Range("E7").Value = 1.55

$ws.Range("D8").Value = "Exit 3 Ratio"
$ws.Range("D9").Value = "Exit 3 Price"
$ws.Range("D10").Value = "Exit 4 Ratio"
$ws.Range("D11").Value = "Exit 4 Price"

# --- "Result" header + formula (D1:E1) ---
$ws.Range("D1").Value = "Result"
$ws.Range("E1").Formula = "=`$H`$3*`$E`$4 +`$H`$4*`$E`$6+`$H`$5*`$E`$8+`$H`$6*E10"

# --- Exit R-multiple helper labels/formulas (G column), entered in
#     G4, G3, G5, G6 order to match original authoring order. ---
$ws.Range("G4").Value = "Exit2 R"
$ws.Range("H4").Formula = "=IF(ISBLANK(E7), 0, (E7-`$E`$2) / `$H`$2)"

$ws.Range("G3").Value = "Exit1 R"
$ws.Range("H3").Formula = "=IF(ISBLANK(E5), 0, (E5-`$E`$2) / `$H`$2)"

$ws.Range("G5").Value = "Exit3 R"
$ws.Range("H5").Formula = "=IF(ISBLANK(E9), 0, (E9-`$E`$2) / `$H`$2)"

$ws.Range("G6").Value = "Exit4 R"
$ws.Range("H6").Formula = "=IF(ISBLANK(E11), 0, (E11-`$E`$2) / `$H`$2)"

# --- Mirror of Enter Price / Stop Price (D2:E3) + R (pips) helper (G2:H2).
#     These reuse pre-existing shared strings (indices 0/1/2) so their
#     position doesn't affect the new-string ordering above. ---
$ws.Range("D2").Value = "Enter Price"
$ws.Range("E2").Formula = "=B2"

$ws.Range("D3").Value = "Stop Price"
$ws.Range("E3").Formula = "=B3"

$ws.Range("G2").Value = "R (pips)"
$ws.Range("H2").Formula = "=E2-E3"

# --- Column D width to match column A, and new selection anchor ---
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Range("D1:E7").Select()

$wb.Application.Calculate()
